$wb = $excel.ActiveWorkbook

# --- Sheet1 ("SNUC5_noCTRL_meas", the raw/meas sheet) ---
# This sheet carried a stray leftover index column (row numbers only, no
# other data) all the way down to row 87, even though the real data block
# only spans rows 1:44. Trim the sheet back down to the real data by
# deleting those extra trailing rows.
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate()
$ws1.Range("A45:A87").EntireRow.Delete()

# Sheet1 becomes the sheet the workbook opens on (previously Sheet3 was
# active), with the view left scrolled/selected where the author was last
# working.
$ws1.Range("D56").Select()
